$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$aw = $excel.ActiveWindow
$panes = $aw.Panes
$p2 = $panes.Item(2)
$p2.ScrollRow = 59
$p2.ScrollColumn = 1
$ws.Range("C81").Select() | Out-Null
